# Grid_Wire_Channel_Mapping.xlsx -- "Fixed optimized code bug"
#
# The channel-remap table in columns B (Wires 20x4 "Channel") and F
# (Wires 16x4 "Channel") was recomputed. Column B gets entirely new
# literal values for rows 3-82; column F gets new literal values for
# most rows, but rows 24-49 keep (an updated) shared formula
# "=F8+16" filled down, anchored at F24 (rows 25-28 are subsequently
# overwritten with literal numbers, same as in the target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Column B ("Channel" for Wires 20x4), rows 3-82
# ---------------------------------------------------------------
$bValues = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,76,77,78,79,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,72,73,74,75,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,68,69,70,71,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $ws.Cells.Item(3 + $i, 2).Value = $bValues[$i]
}

# ---------------------------------------------------------------
# Column F ("Channel" for Wires 16x4), rows 3-23: plain values
# ---------------------------------------------------------------
$fTop = @(14,15,12,13,10,11,8,9,6,7,4,5,2,3,0,1,30,31,28,29,26)
for ($i = 0; $i -lt $fTop.Length; $i++) {
    $ws.Cells.Item(3 + $i, 6).Value = $fTop[$i]
}

# ---------------------------------------------------------------
# Column F, rows 24-49: shared formula "=F8+16" filled down
# (matches master cell landing on F24 with ref F24:F49)
# ---------------------------------------------------------------
$ws.Range("F24:F49").Formula = "=F8+16"

# Rows 25-28 get overwritten with literal values afterwards (the
# formula fill above is only "good" for F24 and F29:F49).
$fHole = @(24,25,22,23)
for ($i = 0; $i -lt $fHole.Length; $i++) {
    $ws.Cells.Item(25 + $i, 6).Value = $fHole[$i]
}

# ---------------------------------------------------------------
# Column F, rows 50-66: plain values
# ---------------------------------------------------------------
$fBottom = @(33,62,63,60,61,58,59,56,57,54,55,52,53,50,51,48,49)
for ($i = 0; $i -lt $fBottom.Length; $i++) {
    $ws.Cells.Item(50 + $i, 6).Value = $fBottom[$i]
}

# ---------------------------------------------------------------
# View state: selection moved to F66, scrolled back to A1 (no
# frozen/forced topLeftCell anymore).
# ---------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("F66").Select()
